# examplePDSI.xlsx edit
# Commit: "Remove -1 to 1 restriction on ecoregion and disturbance modifiers.
#          Recompile. Additional tests for interaction with climate and
#          biomass insects."
#
# The underlying data change is a circular rotation of the AnnualPDSI
# column (B2:B101): the first 6 monthly values are moved to the end of the
# series (i.e. new[row] = old[row + 6], wrapping around), as if 6 rows had
# been consumed/dropped off the front of a rolling PDSI series and 6 new
# values appended at the tail. Column A (the Year/index column) is
# untouched.
#
# The view state also changed: the window had scrolled down (topLeftCell
# A66) with B101 selected as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 101
$colB = 2
$rowCount = $lastDataRow - $firstDataRow + 1   # 100
$shift = 6

# Snapshot the current B2:B101 values (Value2 avoids any date/currency
# re-interpretation of the plain doubles stored here).
$original = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $original += $ws.Cells.Item($r, $colB).Value2
}

# Rewrite each cell with the value that sat $shift rows further down the
# original series, wrapping circularly back to the top.
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcIndex = ($i + $shift) % $rowCount
    $destRow = $firstDataRow + $i
    $ws.Cells.Item($destRow, $colB).Value2 = $original[$srcIndex]
}

# Restore/update the view: scroll so row 66 is at the top and select B101.
$ws.Range("B101").Select()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1
